$wb = $excel.ActiveWorkbook

# graphSpec sheet holds the plot configuration table
$ws = $wb.Worksheets.Item("graphSpec")

# Rename the "Condition_Var" header to "plotEnabled" and optimize the
# condTrue placeholder values into a real boolean TRUE flag.
$ws.Range("G2").Value = "plotEnabled"

$ws.Range("G4").Value = $true
$ws.Range("G6").Value = $true
$ws.Range("G7").Value = $true
$ws.Range("G8").Value = $true
$ws.Range("G9").Value = $true
$ws.Range("G10").Value = $true

$ws.Range("G5").Select()
